# --- GEVO yearly financials update ---
# A new reporting period (newest year) is inserted as column D,
# pushing the existing D:K figures one column to the right (E:L).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()

# Carry the number formatting (date / #,##0) from column E into the
# freshly inserted column D for each of the three statement blocks
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the newest periods figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 32900
$ws.Range("D9").Value = 41600
$ws.Range("D10").Value = -8700
$ws.Range("D12").Value = 5400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 2200
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 57300
$ws.Range("D18").Value = -24400
$ws.Range("D20").Value = -300
$ws.Range("D21").Value = -18100
$ws.Range("D22").Value = 3200
$ws.Range("D23").Value = -28000
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -28000
$ws.Range("D27").Value = -28000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 300
$ws.Range("D33").Value = -28000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -28000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 33700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 500
$ws.Range("D44").Value = 3200
$ws.Range("D45").Value = 1300
$ws.Range("D46").Value = 38700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 67000
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 107000
$ws.Range("D57").Value = 4900
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 400
$ws.Range("D60").Value = 5300
$ws.Range("D61").Value = 12600
$ws.Range("D62").Value = 400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 18200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -429300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 88800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -28000
$ws.Range("D83").Value = 6500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -15900
$ws.Range("D91").Value = -2200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 40300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 22200

# A handful of cells in the shifted range also received updated figures
# (not a pure shift of the prior values) -- apply those corrections explicitly
$ws.Range("E21").Value = -14900
$ws.Range("F21").Value = -22800
$ws.Range("G21").Value = -23100
$ws.Range("H21").Value = -29500
$ws.Range("I21").Value = -54200
$ws.Range("J21").Value = "NA"
$ws.Range("K21").Value = -40000
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 500
$ws.Range("E100").Value = 3600
$ws.Range("F100").Value = 37300
$ws.Range("G100").Value = 40200
$ws.Range("H100").Value = 28200
$ws.Range("I100").Value = 12600
$ws.Range("J100").Value = 93600
$ws.Range("K100").Value = 120600
$ws.Range("E102").Value = -18900
$ws.Range("F102").Value = 10900
$ws.Range("G102").Value = 10700
$ws.Range("H102").Value = -18300
$ws.Range("I102").Value = -42100
$ws.Range("J102").Value = -27500
$ws.Range("K102").Value = 79000
